$p = $ppt.ActivePresentation

# --- 1) Re-style the three data tables (slides 14, 15, 16) with the new
#     built-in table style GUID (was {D3283E9F-2BB7-47D1-AAB4-0AB70200B631},
#     now {CBD910D5-2148-4AAF-A9EE-1538515ECB2C}). Table styles can't be
#     assigned through the Style/StyleId properties, they need ApplyStyle().
$newStyleId = "{CBD910D5-2148-4AAF-A9EE-1538515ECB2C}"
foreach ($slideIdx in 14,15,16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newStyleId)
        }
    }
}

# --- 2) Swap the presentation's colour theme from "Integral" (Red Violet)
#     back to the stock "Office Theme" palette. The deck's applied design
#     is reachable through the Slide ThemeColorScheme (12 theme colour
#     slots, same order as the OOXML clrScheme: dk1, lt1, dk2, lt2,
#     accent1-6, hlink, folHlink).
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $tcs.Item($i).RGB = $officeColors[$i - 1]
}
